$d = $word.ActiveDocument

function Replace-ExactText($old, $new) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Find failed for: $old"
    }
    $r.Text = $new
}

Replace-ExactText "{Loving Kindness Meditation} " "{Méditation de l'Amour Bienveillant} "
Replace-ExactText "Today's pause is called loving kindness meditation." "La pause du jour s'appelle la méditation de l'amour bienveillant."
Replace-ExactText "Find a comfortable sitting position, your feet flat on the floor, your hands resting in your lap." "Asseyez-vous confortablement, les pieds à plat sur le sol, les mains posées sur les genoux."
Replace-ExactText "Close your eyes if you are comfortable, or maintain a soft gaze. " "Fermez les yeux si vous êtes à l'aise, ou maintenez un regard doux. "
Replace-ExactText "Ask yourself, “What is my experience in this moment?” " "Demandez-vous, “Quelle est mon expérience en ce moment précis?” "
Replace-ExactText "Notice what thoughts you are experiencing. " "Prêtez attention aux pensées qui vous viennent. "
Replace-ExactText "Notice how you feel emotionally. " "Prêtez attention aux émotions que vous ressentez. "
Replace-ExactText "Remarquez ce que vous ressentez dans votre corps. Notice any discomfort or tension." "Remarquez ce que vous ressentez dans votre corps. Prêtez attention à toute gêne ou tension que vous pourrez ressentir."
Replace-ExactText "Connect to your heart in a kind and gentle way. You may want to place one hand on your heart or chest. " "Connectez-vous à votre cœur de manière aimable et douce. Vous pouvez placer une main sur votre cœur ou poitrine. "
Replace-ExactText "You can then say the following words silently to yourself: " "Vous pouvez ensuite vous dire silencieusement les mots suivants : "
Replace-ExactText "May I be peaceful. " "Que je sois en paix. "
Replace-ExactText "May I be safe. " "Que je sois en sécurité. "
Replace-ExactText "May I be healthy. " "Que je sois en bonne santé. "
Replace-ExactText "May I be happy. " "Que je sois heureux. "
Replace-ExactText "May I feel loved. " "Que je sois aimé. "
Replace-ExactText "Repeat slowly once or twice, taking your time between each phrase." "Répétez lentement une ou deux fois, en prenant votre temps entre chaque phrase."
Replace-ExactText "If you feel comfortable, you can also send thoughts of loving-kindness to your child, your partner, your family, and anyone else who is close to you in your life." "Si vous vous sentez à l'aise, vous pouvez aussi envoyer des pensées d'amour bienveillant à votre enfant, votre partenaire, votre famille, et toute autre personne qui vous est proche."
Replace-ExactText "Now, allow your focus to expand to the whole body. " "Maintenant, laissez votre attention s'étendre à tout votre corps. "
Replace-ExactText "Allow your focus to expand to the sounds in the room. " "Laissez votre attention s'étendre aux sons dans la pièce dans laquelle vous vous trouvez. "
Replace-ExactText "Open your eyes, and when you are ready, continue to your lesson with a sense of calm. " "Ouvrez les yeux, et lorsque vous êtes prêt, continuez votre leçon en toute sérénité. "
Replace-ExactText "Thank you for taking a moment to pause with us. " "Merci d'avoir pris le temps de faire une pause avec nous. "
